$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new test case row value in column A, row 2 ("run")
$ws.Range("A2").Value = "run"

# Update the active selection shown in the sheet view
$ws.Range("C8").Select()
